{"js": "// Replace the date line and the multiplication problems, matching the\n// unified diff exactly. Each \"before\" text is unique in the document, so\n// a search-and-replace keyed on old text is unambiguous and preserves the\n// run formatting (rFonts / sz) since we replace in-place inside the found\n// range.\nconst replacements = [\n  [\"2024-08-02 Friday\", \"2024-08-03 Saturday\"],\n  [\"511\u00d77=\", \"105\u00d74=\"],\n  [\"135\u00d78=\", \"519\u00d74=\"],\n  [\"931\u00d75=\", \"352\u00d73=\"],\n  [\"206\u00d73=\", \"152\u00d76=\"],\n  [\"110\u00d74=\", \"844\u00d72=\"],\n  [\"754\u00d76=\", \"615\u00d73=\"],\n  [\"849\u00d73=\", \"352\u00d73=\"],\n  [\"341\u00d77=\", \"720\u00d78=\"],\n  [\"370\u00d76=\", \"747\u00d74=\"],\n  [\"951\u00d78=\", \"611\u00d74=\"],\n  [\"752\u00d75=\", \"567\u00d73=\"],\n  [\"246\u00d76=\", \"818\u00d78=\"],\n  [\"404\u00d78=\", \"232\u00d72=\"],\n  [\"571\u00d79=\", \"877\u00d78=\"],\n  [\"417\u00d72=\", \"198\u00d77=\"],\n  [\"920\u00d74=\", \"677\u00d74=\"],\n  [\"885\u00d72=\", \"370\u00d73=\"],\n  [\"166\u00d74=\", \"220\u00d72=\"],\n  [\"453\u00d75=\", \"493\u00d72=\"],\n  [\"304\u00d72=\", \"162\u00d73=\"],\n  [\"173\u00d79=\", \"402\u00d75=\"],\n  [\"221\u00d77=\", \"610\u00d78=\"],\n  [\"900\u00d76=\", \"755\u00d78=\"],\n  [\"765\u00d73=\", \"434\u00d78=\"],\n  [\"194\u00d73=\", \"374\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Only the first match is replaced \u2014 every \"before\" string above occurs\n  // exactly once in the source document.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the multiplication problems, matching the\n# unified diff exactly. Each \"before\" text is unique in the document, so\n# Find/Replace keyed on the old text is unambiguous and preserves the run\n# formatting (rFonts / sz) since Word replaces text in place.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-08-02 Friday\", \"2024-08-03 Saturday\"),\n    @(\"511\u00d77=\", \"105\u00d74=\"),\n    @(\"135\u00d78=\", \"519\u00d74=\"),\n    @(\"931\u00d75=\", \"352\u00d73=\"),\n    @(\"206\u00d73=\", \"152\u00d76=\"),\n    @(\"110\u00d74=\", \"844\u00d72=\"),\n    @(\"754\u00d76=\", \"615\u00d73=\"),\n    @(\"849\u00d73=\", \"352\u00d73=\"),\n    @(\"341\u00d77=\", \"720\u00d78=\"),\n    @(\"370\u00d76=\", \"747\u00d74=\"),\n    @(\"951\u00d78=\", \"611\u00d74=\"),\n    @(\"752\u00d75=\", \"567\u00d73=\"),\n    @(\"246\u00d76=\", \"818\u00d78=\"),\n    @(\"404\u00d78=\", \"232\u00d72=\"),\n    @(\"571\u00d79=\", \"877\u00d78=\"),\n    @(\"417\u00d72=\", \"198\u00d77=\"),\n    @(\"920\u00d74=\", \"677\u00d74=\"),\n    @(\"885\u00d72=\", \"370\u00d73=\"),\n    @(\"166\u00d74=\", \"220\u00d72=\"),\n    @(\"453\u00d75=\", \"493\u00d72=\"),\n    @(\"304\u00d72=\", \"162\u00d73=\"),\n    @(\"173\u00d79=\", \"402\u00d75=\"),\n    @(\"221\u00d77=\", \"610\u00d78=\"),\n    @(\"900\u00d76=\", \"755\u00d78=\"),\n    @(\"765\u00d73=\", \"434\u00d78=\"),\n    @(\"194\u00d73=\", \"374\u00d75=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceOne = 2\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
